$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "tuned algorithm": every division now plays the same number of max games (50)
# and the minimum days between games for each division is reduced to 1.
$ws.Range("F2:F8").Value = 50
$ws.Range("H2:H8").Value = 1

# Leave the cursor where the author last left it when saving.
$ws.Range("H9").Select()
